# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# Price cells in column D that look purely numeric ("211.10", "7.01", ...)
# are forced to Text format before assignment so Excel keeps them as the
# exact literal string (preserving trailing zeros / not coercing to a
# float), then the cell style is reset back to "Normal" so no stray
# number-format style lingers on the cell afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.618.97'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.596.21'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0838'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '1.820.07'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '1.607.98'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = '26.600.52'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0510'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").Value = '1.276.93'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.617'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.75%  '
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +18.07%  '
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("D45").Value = '1.732.57'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("E48").Value = '  +3.77%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.39%  '
